$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) data cells to Text format so numeric-looking
# strings like "1.002" are preserved exactly as text, matching source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.222.88'
$ws.Range("E2").Value = '  -0.23%  '

$ws.Range("D3").Value = '1.870.05'
$ws.Range("E3").Value = '  +3.28%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '311.71'
$ws.Range("E5").Value = '  -0.25%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("E7").Value = '  -2.19%  '

$ws.Range("D8").Value = '0.3951'
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("D9").Value = '0.09922'
$ws.Range("E9").Value = '  +26.90%  '

$ws.Range("D10").Value = '1.129'
$ws.Range("E10").Value = '  +1.82%  '

$ws.Range("D11").Value = '41.16'
$ws.Range("E11").Value = '  +0.26%  '

$ws.Range("D12").Value = '6.468'
$ws.Range("E12").Value = '  +1.61%  '

$ws.Range("D13").Value = '20.94'
$ws.Range("E13").Value = '  +2.15%  '

$ws.Range("D14").Value = '1.870.47'
$ws.Range("E14").Value = '  +4.05%  '

$ws.Range("D15").Value = '1.002'
$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("D16").Value = '7.395'
$ws.Range("E16").Value = '  +0.74%  '

$ws.Range("D17").Value = '0.00001138'
$ws.Range("E17").Value = '  +5.46%  '

$ws.Range("D18").Value = '93.75'
$ws.Range("E18").Value = '  +1.18%  '

$ws.Range("D19").Value = '0.06652'
$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '17.44'
$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").Value = '6.114'
$ws.Range("E22").Value = '  +1.63%  '

$ws.Range("D23").Value = '28.283.01'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").Value = '11.35'
$ws.Range("E24").Value = '  +1.89%  '

$ws.Range("D25").Value = '2.262'
$ws.Range("E25").Value = '  +1.26%  '

$ws.Range("D26").Value = '2.550'
$ws.Range("E26").Value = '  +3.83%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '21.28'
$ws.Range("E27").Value = '  +3.82%  '

$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.085.32'
$ws.Range("E28").Value = '  +3.70%  '

$ws.Range("D29").Value = '157.73'
$ws.Range("E29").Value = '  -1.95%  '

$ws.Range("D30").Value = '127.81'
$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("D31").Value = '0.1059'
$ws.Range("E31").Value = '  -3.22%  '

$ws.Range("D32").Value = '1.057'
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("D33").Value = '5.625'
$ws.Range("E33").Value = '  +0.98%  '

$ws.Range("D34").Value = '3.603'
$ws.Range("E34").Value = '  -1.37%  '

$ws.Range("D35").Value = '0.06812'
$ws.Range("E35").Value = '  -4.70%  '

$ws.Range("D36").Value = '9.446'
$ws.Range("E36").Value = '  +3.14%  '

$ws.Range("D37").Value = '0.02396'
$ws.Range("E37").Value = '  +1.86%  '

$ws.Range("D38").Value = '0.2187'
$ws.Range("E38").Value = '  +0.46%  '

$ws.Range("D39").Value = '5.016'
$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("E40").Value = '  -0.93%  '

$ws.Range("D41").Value = '0.6302'
$ws.Range("E41").Value = '  +1.95%  '

$ws.Range("D42").Value = '1.177'
$ws.Range("E42").Value = '  +1.78%  '

$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").Value = '13.55'
$ws.Range("E44").Value = '  +2.47%  '

$ws.Range("D45").Value = '0.6017'
$ws.Range("E45").Value = '  +0.84%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.667'
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '1.274'
$ws.Range("E47").Value = '  -2.31%  '

$ws.Range("D48").Value = '124.86'
$ws.Range("E48").Value = '  -0.41%  '

$ws.Range("D49").Value = '1.994'
$ws.Range("E49").Value = '  +3.60%  '

$ws.Range("D50").Value = '1.201'
$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '1.124'
$ws.Range("E51").Value = '  +5.30%  '
